$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (Excel serial 45206 = 2023-10-07)
# that must be updated to serial 45208 (2023-10-09) for every data row
# (rows 2 through 122).
for ($row = 2; $row -le 122; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
